$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 766206
$ws.Range("E2").Value = 1428873872

$ws.Range("C10").Value = 345448
$ws.Range("E10").Value = 1816668803

$ws.Range("C13").Value = 187760
$ws.Range("E13").Value = 1162504460

$ws.Range("C36").Value = 211189
$ws.Range("E36").Value = 404207167

$ws.Range("C57").Value = 31589
$ws.Range("E57").Value = 162167845

$ws.Range("C67").Value = 27091
$ws.Range("E67").Value = 168553316

$ws.Range("C72").Value = 331301
$ws.Range("E72").Value = 635348916

$ws.Range("C78").Value = 178404
$ws.Range("E78").Value = 892005395

$ws.Range("C79").Value = 679
$ws.Range("E79").Value = 20307920

$ws.Range("C91").Value = 18404
$ws.Range("E91").Value = 72097061

$ws.Range("C93").Value = 16554
$ws.Range("E93").Value = 48244387

$ws.Range("C112").Value = 145185
$ws.Range("E112").Value = 715559271

$ws.Range("C115").Value = 81783
$ws.Range("D115").Value = 14447
$ws.Range("E115").Value = 435968262

$ws.Range("C121").Value = 1305802
$ws.Range("E121").Value = 2273392960

$ws.Range("C127").Value = 9136
$ws.Range("E127").Value = 110233894

$ws.Range("C128").Value = 279
$ws.Range("E128").Value = 5709119

$ws.Range("C129").Value = 632748
$ws.Range("E129").Value = 3417175036

$ws.Range("C130").Value = 4226
$ws.Range("E130").Value = 139004203

$ws.Range("C132").Value = 585033
$ws.Range("E132").Value = 3441706448

$ws.Range("C136").Value = 26629
$ws.Range("E136").Value = 141829041

$ws.Range("C144").Value = 24512
$ws.Range("E144").Value = 88163649

$ws.Range("C151").Value = 39255
$ws.Range("E151").Value = 59799025

$ws.Range("C154").Value = 17967
$ws.Range("E154").Value = 69446902

$ws.Range("C157").Value = 629
$ws.Range("E157").Value = 1385587

$ws.Range("C171").Value = 95809
$ws.Range("E171").Value = 490329269

$ws.Range("C196").Value = 595463
$ws.Range("E196").Value = 983954243

$ws.Range("C215").Value = 230234
$ws.Range("E215").Value = 408661646

$ws.Range("C221").Value = 135466
$ws.Range("E221").Value = 681712088

$ws.Range("C229").Value = 612508
$ws.Range("E229").Value = 1040624400

$ws.Range("C237").Value = 283247
$ws.Range("E237").Value = 1437565257

$ws.Range("C240").Value = 205847
$ws.Range("E240").Value = 1066946940
